$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update URL, Version, Date, Publisher ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/procedure-group"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: clear the Constraint(s) note on the root Extension row,
#     and update the payer-procedure-group value set URL ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/payer-procedure-group"

# The longer URL text widens the auto-fit "Binding Value Set" column (Y / col 25)
$elements.Columns.Item(25).ColumnWidth = 61.165
